$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1997,92534949.92261896,366107152.671737,603639380.3493997),
    @(1998,73015573.13914414,144842530.7588122,694144378.455973),
    @(1999,127533577.158244,94505985.09341794,641294859.4167044),
    @(2000,117717354.2189644,85259607.93777138,765801677.4468803),
    @(2001,193613346.4684302,55877348.80975366,870977414.0426997),
    @(2002,179518694.6202532,902645898.4177215,1078404236.075949),
    @(2003,213193731.4856761,939701385.2872183,982696186.8700914),
    @(2004,208492822.2225001,563701390.6239058,1441165959.17788),
    @(2005,153970883.8277984,511558922.3607715,1783372958.611462),
    @(2006,155630768.5524178,499792471.9320854,1485635484.890729),
    @(2007,77777261.54150286,621694178.2528318,1540099115.439557),
    @(2008,95870212.88534531,626151942.0153413,1627125898.616248),
    @(2009,96309323.66573989,657129917.6624393,1657641692.098827),
    @(2010,87563098.32237658,764624986.1314064,1993118808.239646),
    @(2011,82492684.56839997,785255510.4660183,2207040594.443008),
    @(2012,88806678.39962138,940284716.7317612,2194751836.958079),
    @(2013,231465122.7407945,1068800287.814614,2130230479.616107),
    @(2014,266703630.3963684,1086533080.168774,2098966207.998798),
    @(2015,231359992.4266163,1015386024.421622,1921532195.557813),
    @(2016,206234994.1223276,924372641.5920178,1882075451.245157),
    @(2017,207922967.828885,1097275594.047953,1973059458.032405),
    @(2018,216621050.4713161,1051335925.556159,1919784570.199054),
    @(2019,247194582.859546,1103886780.764249,1881913482.410036),
    @(2020,275570037.3710343,919798135.8751496,1901960179.090258),
    @(2021,325882679.0956431,958778950.811841,2110865612.403105),
    @(2022,320859386.8942198,947217707.0477858,2255243528.627884),
    @(2023,342904604.6797927,948722283.3989342,2616138411.594326),
    @(2024,321365149.12,908436601.3599999,2884059718.07)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
